# Apply "new basic stats 2022" update to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (columns A: label, B: freq, C: pct, D: var name stays the same)
$varName = "all_museums-subject_matter_simpl_aggr"

$data = @(
    @{ Label = "SMALL_SUBJECTS";            Freq = 1626; Pct = 38.45 },
    @{ Label = "NOT_AVAIL";                 Freq = 1000; Pct = 23.65 },
    @{ Label = "houses";                    Freq = 517;  Pct = 12.23 },
    @{ Label = "other";                     Freq = 361;  Pct = 8.539999999999999 },
    @{ Label = "fine_and_decorative_arts";  Freq = 196;  Pct = 4.63 },
    @{ Label = "trains_and_railways";       Freq = 143;  Pct = 3.38 },
    @{ Label = "mixed";                     Freq = 136;  Pct = 3.22 },
    @{ Label = "regiment";                  Freq = 131;  Pct = 3.1 },
    @{ Label = "encyclopaedic";             Freq = 119;  Pct = 2.81 }
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item.Label
    $ws.Cells.Item($row, 2).Value = $item.Freq
    $ws.Cells.Item($row, 3).Value = $item.Pct
    $ws.Cells.Item($row, 4).Value = $varName
    $row++
}

# Remove old rows 11-15 (which existed previously but are no longer part of the data)
$ws.Range("A11:D15").ClearContents()
